# "Adapt tests to control version"
# Add a "version" column to the GeoJSON form's "settings" sheet:
#   C1 = "version" (new header, default/unstyled like the rest of row 2)
#   C2 = 1         (numeric version value)
# and leave the selection on the newly added cell below the data (C3),
# matching the end-user edit captured in the workbook diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("settings")

$ws.Range("C1").Value = "version"
$ws.Range("C2").Value = 1

$ws.Range("C3").Select() | Out-Null
